# Apply dataset tracking edits to the workbook
$wb = $excel.ActiveWorkbook

# --- Sheet: Dataset Registry ---
$ws1 = $wb.Worksheets.Item("Dataset Registry")
$ws1.Range("B2").Value = "config-1.yml"
$ws1.Range("D2").NumberFormat = "@"
$ws1.Range("D2").Value = "2025-07-03"
$ws1.Range("B3").Value = "config-2.yml"
$ws1.Range("C3").Value = "n1000000_f_init20_cont20_disc0_add0_pert-none_scl0_func-linear_noise0_dataset.csv"
$ws1.Range("D3").NumberFormat = "@"
$ws1.Range("D3").Value = "2025-07-03"

# --- Sheet: Configuration Details ---
$ws2 = $wb.Worksheets.Item("Configuration Details")
$ws2.Range("B2").Value = "config-1.yml"
$ws2.Range("B3").Value = "config-2.yml"
$ws2.Range("F3").Value = 20
$ws2.Range("G3").Value = 0

# --- Sheet: Feature Details ---
$ws3 = $wb.Worksheets.Item("Feature Details")
for ($r = 22; $r -le 41; $r++) {
    $ws3.Range("C$r").Value = "continuous"
}

# --- Sheet: File Metadata ---
$ws5 = $wb.Worksheets.Item("File Metadata")
$ws5.Range("B2").Value = "configs\data_generation\config-1.yml"
$ws5.Range("B3").Value = "configs\data_generation\config-2.yml"
$ws5.Range("C3").Value = "data\n1000000_f_init20_cont20_disc0_add0_pert-none_scl0_func-linear_noise0_dataset.csv"
$ws5.Range("D3").Value = "reports\figures\n1000000_f_init20_cont20_disc0_add0_pert-none_scl0_func-linear_noise0_plot.pdf"
$ws5.Range("E3").Value = 379.1
$ws5.Range("F3").Value = "04111418c3"
$ws5.Range("G3").Value = "Linear function, 20C/0D features"
